$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Sprint" column header in G1
$ws.Range("G1").Value = "Sprint"

# Fill G2:G17 with "S1"
$ws.Range("G2:G17").Value = "S1"

# Update selection to match the target state
$ws.Range("G2:G17").Select()
